# update enzyme modifying therapy
#
# Appends four new pattern rows (50-53) to Sheet1:
#   - gout suppressant            -> antigout therapy
#   - antirheumatic drug          -> antirheumatic agent therapy
#   - acetylcholinesterase inhibitor -> EC 3.1.1.7 (acetylcholinesterase) inhibitor  (highlighted)
#   - cholinesteriase inhibitor   -> EC 3.1.1.8 (cholinesterase) inhibitor (highlighted)
#
# Columns are: A=defined_class, B=defined_class_name, C=chemical_identity, D=chemical_identity_label

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 50 : gout suppressant -----------------------------------------
$ws.Range("A50").Value = "MAXO_0000164"
$ws.Range("B50").Value = "antigout therapy"
$ws.Range("C50").Value = "CHEBI_35845"
$ws.Range("D50").Value = "gout suppressant"

# ---- Row 51 : antirheumatic drug ---------------------------------------
$ws.Range("A51").Value = "MAXO_0000643"
$ws.Range("B51").Value = "antirheumatic agent therapy"
$ws.Range("C51").Value = "CHEBI_35842"
$ws.Range("D51").Value = "antirheumatic drug"

# ---- Row 52 : acetylcholinesterase inhibitor (highlighted row) --------
# C52 re-uses the same highlighted "definition" style already used elsewhere (e.g. C45)
$ws.Range("C45").Copy($ws.Range("C52")) | Out-Null
$ws.Range("C52").Value = "CHEBI:38462"

$ws.Range("B52").Value = "acetylcholinesterase inhibitor"

# A52 gets a brand-new highlight font: 18pt Helvetica, teal (RGB 0,124,130)
$ws.Range("C44").Copy($ws.Range("A52")) | Out-Null
$ws.Range("A52").Value = "MAXO_0000645"
$hiFont = $ws.Range("A52").Font
$hiFont.Size = 18
$hiFont.Color = 8551424

# D52 shares the exact same highlight style as A52
$ws.Range("A52").Copy($ws.Range("D52")) | Out-Null
$ws.Range("D52").Value = "EC 3.1.1.7 (acetylcholinesterase) inhibitor"

$ws.Rows.Item(52).RowHeight = 23

# ---- Row 53 : cholinesterase inhibitor ---------------------------------
$ws.Range("A53").Value = "MAXO_0000210"
$ws.Range("B53").Value = "cholinesteriase inhibitor"
$ws.Range("C53").Value = "CHEBI:37733"

# D53 shares the same highlight style as A52/D52
$ws.Range("A52").Copy($ws.Range("D53")) | Out-Null
$ws.Range("D53").Value = "EC 3.1.1.8 (cholinesterase) inhibitor"

$ws.Rows.Item(53).RowHeight = 23

# ---- mirror the source workbook's final selection ----------------------
$ws.Range("D53").Select() | Out-Null

$wb.Save()
